$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value is a "clean" decimal number (single dot).
# Excel's automatic type-detection would silently convert these from the
# original text("inlineStr")-typed cells into numeric cells, which the
# source workbook does not do. Force them to stay text by pre-formatting
# the cell as Text ("@") before writing the value, matching how the
# original file stores these prices as strings.
$textPriceCells = @(5,6,7,10,13,14,16,19,20,22,24,26,27,29,33,34,35,36,40,42,45,46,49,50)
foreach ($r in $textPriceCells) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# --- Row 2 ---
$ws.Range("D2").Value = "37.768.77"
$ws.Range("E2").Value = "  -0.17%  "

# --- Row 3 ---
$ws.Range("D3").Value = "2.039.50"
$ws.Range("E3").Value = "  +0.12%  "

# --- Row 4 ---
$ws.Range("E4").Value = "  +0.01%  "

# --- Row 5 ---
$ws.Range("D5").Value = "227.34"
$ws.Range("E5").Value = "  -0.18%  "

# --- Row 6 ---
$ws.Range("D6").Value = "0.607"
$ws.Range("E6").Value = "  -1.14%  "

# --- Row 7 ---
$ws.Range("D7").Value = "59.49"
$ws.Range("E7").Value = "  -1.26%  "

# --- Row 8 ---
$ws.Range("E8").Value = "  +0.05%  "

# --- Row 9 ---
$ws.Range("E9").Value = "  -2.83%  "

# --- Row 10 ---
$ws.Range("D10").Value = "0.0837"
$ws.Range("E10").Value = "  +2.54%  "

# --- Row 11 ---
$ws.Range("E11").Value = "  -0.33%  "

# --- Row 12 ---
$ws.Range("D12").Value = "2.342.93"
$ws.Range("E12").Value = "  +0.16%  "

# --- Row 13 ---
$ws.Range("D13").Value = "14.41"
$ws.Range("E13").Value = "  -1.72%  "

# --- Row 14 ---
$ws.Range("D14").Value = "21.01"
$ws.Range("E14").Value = "  -0.43%  "

# --- Row 15 ---
$ws.Range("E15").Value = "  +4.79%  "

# --- Row 16 ---
$ws.Range("D16").Value = "0.769"
$ws.Range("E16").Value = "  +0.84%  "

# --- Row 17 ---
$ws.Range("D17").Value = "2.043.50"
$ws.Range("E17").Value = "  +0.03%  "

# --- Row 18 ---
$ws.Range("D18").Value = "37.749.60"
$ws.Range("E18").Value = "  -0.24%  "

# --- Row 19 ---
$ws.Range("D19").Value = "69.45"
$ws.Range("E19").Value = "  -0.65%  "

# --- Row 20 ---
$ws.Range("D20").Value = "5.91"
$ws.Range("E20").Value = "  -2.80%  "

# --- Row 21 ---
$ws.Range("D21").Value = "0.0₃0823"
$ws.Range("E21").Value = "  -0.20%  "

# --- Row 22 ---
$ws.Range("D22").Value = "223.71"
$ws.Range("E22").Value = "  -0.76%  "

# --- Row 23 ---
$ws.Range("E23").Value = "  +0.55%  "

# --- Row 24 ---
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  +1.42%  "

# --- Row 25 ---
$ws.Range("E25").Value = "  +2.73%  "

# --- Row 26 ---
$ws.Range("D26").Value = "169.55"
$ws.Range("E26").Value = "  +2.54%  "

# --- Row 27 ---
$ws.Range("D27").Value = "9.35"
$ws.Range("E27").Value = "  +0.71%  "

# --- Row 28 ---
$ws.Range("E28").Value = "  -0.74%  "

# --- Row 29 ---
$ws.Range("D29").Value = "18.78"
$ws.Range("E29").Value = "  -0.83%  "

# --- Row 30 ---
$ws.Range("E30").Value = "  -0.15%  "

# --- Row 31 ---
$ws.Range("E31").Value = "  -1.04%  "

# --- Row 32 ---
$ws.Range("E32").Value = "  +9.14%  "

# --- Row 33 ---
$ws.Range("D33").Value = "4.38"
$ws.Range("E33").Value = "  -1.56%  "

# --- Row 34 (was Hedera, becomes InternetComputer(DFINITY)) ---
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "4.48"
$ws.Range("E34").Value = "  -0.28%  "

# --- Row 35 (was InternetComputer(DFINITY), becomes Hedera) ---
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.0601"
$ws.Range("E35").Value = "  -0.38%  "

# --- Row 36 ---
$ws.Range("D36").Value = "6.56"
$ws.Range("E36").Value = "  +1.60%  "

# --- Row 37 ---
$ws.Range("E37").Value = "  +3.47%  "

# --- Row 39 ---
$ws.Range("E39").Value = "  -0.10%  "

# --- Row 40 ---
$ws.Range("D40").Value = "17.98"
$ws.Range("E40").Value = "  +6.47%  "

# --- Row 41 ---
$ws.Range("D41").Value = "1.526.42"
$ws.Range("E41").Value = "  -1.00%  "

# --- Row 42 ---
$ws.Range("D42").Value = "97.41"
$ws.Range("E42").Value = "  +0.59%  "

# --- Row 43 ---
$ws.Range("E43").Value = "  -1.10%  "

# --- Row 44 ---
$ws.Range("E44").Value = "  +0.46%  "

# --- Row 45 ---
$ws.Range("D45").Value = "0.0907"
$ws.Range("E45").Value = "  -1.86%  "

# --- Row 46 ---
$ws.Range("D46").Value = "4.19"
$ws.Range("E46").Value = "  +7.00%  "

# --- Row 47 ---
$ws.Range("E47").Value = "  -0.35%  "

# --- Row 48 ---
$ws.Range("E48").Value = "  -0.08%  "

# --- Row 49 (was MXToken, becomes FraxShare) ---
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "7.10"
$ws.Range("E49").Value = "  -0.51%  "

# --- Row 50 (was FraxShare, becomes MXToken) ---
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "2.94"
$ws.Range("E50").Value = "  -0.86%  "

# --- Row 51 ---
$ws.Range("D51").Value = "2.232.10"
$ws.Range("E51").Value = "  +0.17%  "
